# Update recalculated BALANCE results across the workbook.
# Values below reflect a fresh recomputation of the aircraft balance
# (center of gravity) estimates after more work on JPADCommander.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value  = 60.58379194471799
$ws.Range("C3").Value  = 12.578175260813666
$ws.Range("C4").Value  = 24.22129792720768
$ws.Range("C5").Value  = 0.5567252463403136
$ws.Range("C7").Value  = 43.62176109828424
$ws.Range("C8").Value  = 12.188303882748276
$ws.Range("C9").Value  = 26.942174991215605
$ws.Range("C10").Value = 0.6192644611368907
$ws.Range("C12").Value = 43.62176109828424
$ws.Range("C13").Value = 12.188303882748276
$ws.Range("C14").Value = 26.942174991215605
$ws.Range("C15").Value = 0.6192644611368907
$ws.Range("C17").Value = 54.03205126990302
$ws.Range("C18").Value = 12.427583849399376
$ws.Range("C19").Value = 17.554146471108247
$ws.Range("C20").Value = 0.4034811242482551
$ws.Range("C22").Value = 51.99537149805351
$ws.Range("C23").Value = 12.380770873021856
$ws.Range("C24").Value = 24.75311268810642
$ws.Range("C25").Value = 0.5689489803721794
$ws.Range("C27").Value = 29.539801660785603
$ws.Range("C28").Value = 66.10696498785293

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value  = 1.3305081214942964
$ws.Range("C6").Value  = 12.330508121494294
$ws.Range("C11").Value = 1.3305081214942964

$ws = $wb.Worksheets.Item("FUEL TANK")
$ws.Range("C2").Value = 1.089073330329886
$ws.Range("C6").Value = 12.089073330329885

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value  = 2.5640710746536346
$ws.Range("C6").Value  = 24.164071074653634
$ws.Range("C11").Value = 2.5640710746536346

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C9").Value  = 10.3174
$ws.Range("C11").Value = 1.0289999999999997
$ws.Range("C20").Value = 10.3174
$ws.Range("C22").Value = 1.0289999999999997

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C9").Value  = 9.633999999999999
$ws.Range("C11").Value = 1.0289999999999997
$ws.Range("C20").Value = 9.633999999999999
$ws.Range("C22").Value = 1.0289999999999997

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 12.321363937165607
$ws.Range("C4").Value = -1.874638699878493
